$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.019.79'
$ws.Range('E2').Value = '  +2.81%  '
$ws.Range('D3').Value = '2.468.49'
$ws.Range('E3').Value = '  +4.68%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '565.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '142.47'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +7.14%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  +0.50%  '
$ws.Range('D9').Value = '2.467.65'
$ws.Range('E9').Value = '  +4.80%  '
$ws.Range('E10').Value = '  +2.01%  '
$ws.Range('E11').Value = '  +0.90%  '
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('E13').Value = '  +2.87%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '26.60'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +9.34%  '
$ws.Range('D15').Value = '2.910.91'
$ws.Range('E15').Value = '  +4.78%  '
$ws.Range('D16').Value = '62.895.21'
$ws.Range('E16').Value = '  +2.87%  '
$ws.Range('E17').Value = '  +3.75%  '
$ws.Range('D18').Value = '2.469.86'
$ws.Range('E18').Value = '  +4.82%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.22'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +3.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '340.43'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +7.05%  '
$ws.Range('E21').Value = '  +2.64%  '
$ws.Range('E22').Value = '  +1.41%  '
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '65.60'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.56%  '
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.49'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +4.33%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.10'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.38'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.08%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.83'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +10.72%  '
$ws.Range('E31').Value = '  +5.69%  '
$ws.Range('D32').Value = '0.0₃0797'
$ws.Range('E32').Value = '  +6.84%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '174.80'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.51'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +9.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.399'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.46%  '
$ws.Range('E36').Value = '  +3.36%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '373.84'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.76%  '
$ws.Range('E39').Value = '  +4.00%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.69'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +8.80%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.33'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +5.18%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '149.80'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +6.67%  '
$ws.Range('E44').Value = '  +3.74%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '20.57'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +5.07%  '
$ws.Range('E46').Value = '  +4.42%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0960'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0516'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.20%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0226'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.79%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0232'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('E51').Value = '  +3.47%  '
